$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing D14 comment text (session 7 / Hypothesis row)
$ws.Cells.Item(14, 4).Value = "very dry! Intertwine with R code of w9? (Theory of the theory of R session)"

# New rows 15-22 (sessions 8-15), following the style of the existing table
$newRows = @(
    @{ Row = 15; A = 8;  B = "Modeling Relationships";            C = "no";  D = "Theory session of w9. intertwine with R code?" },
    @{ Row = 16; A = 9;  B = "Statistical Analysis";               C = "yes" },
    @{ Row = 17; A = 10; B = "GLM" },
    @{ Row = 18; A = 11; B = "GLM R" },
    @{ Row = 19; A = 12; B = "LMM" },
    @{ Row = 20; A = 13; B = "LMM R" },
    @{ Row = 21; A = 14; B = "Question about report" },
    @{ Row = 22; A = 15; B = "extra session: reproducibility" }
)

foreach ($r in $newRows) {
    $rowIndex = $r.Row

    $cellA = $ws.Cells.Item($rowIndex, 1)
    $cellA.Value = $r.A
    $cellA.HorizontalAlignment = -4152  # xlRight, matches styling of column A elsewhere

    if ($r.ContainsKey("B")) { $ws.Cells.Item($rowIndex, 2).Value = $r.B }
    if ($r.ContainsKey("C")) { $ws.Cells.Item($rowIndex, 3).Value = $r.C }
    if ($r.ContainsKey("D")) { $ws.Cells.Item($rowIndex, 4).Value = $r.D }
}

# Final active selection ends up on the first empty row below the table
$ws.Range("A22").Select() | Out-Null
